{"js": "const body = context.document.body;\nconst results = body.search(\"Systemet gemmer beskeden i beskedhistorik.\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found\");\n}\n\nconst target = results.items[0];\ntarget.insertText(\"Systemet gemmer beskeden i beskedhistorik i kronologisk orden.\", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$oldText = \"Systemet gemmer beskeden i beskedhistorik.\"\n$newText = \"Systemet gemmer beskeden i beskedhistorik i kronologisk orden.\"\n\n# wdFindContinue (1) = keep searching past the current selection; wdReplaceOne (1) = replace only the first match.\n$find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n"}
